$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Rename the SmartArt root topic on slide 3 from "高维数据分析" to
#    "数据挖掘与机器学习" (updates both the diagram data part and its cached
#    drawing part).
# ---------------------------------------------------------------------------
$oldTopic = "高维数据分析"
$newTopic = "数据挖掘与机器学习"

$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(1)
if ($sh3.HasSmartArt) {
    $nodes3 = $sh3.SmartArt.AllNodes
    $root3 = $nodes3.Item(1)
    if ($root3.TextFrame2.TextRange.Text -eq $oldTopic) {
        $root3.TextFrame2.TextRange.Text = $newTopic
    }
}

# ---------------------------------------------------------------------------
# 2. Refresh the cached "today" date field shown on every slide layout and on
#    the slide master itself, from 2021/12/27 to 2021/12/28.
# ---------------------------------------------------------------------------
$oldDate = "2021/12/27"
$newDate = "2021/12/28"

$masterShapes = $p.SlideMaster.Shapes
for ($j = 1; $j -le $masterShapes.Count; $j++) {
    $shp = $masterShapes.Item($j)
    if ($shp.PlaceholderFormat.Type -eq 16) {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layoutShapes = $layouts.Item($i).Shapes
    for ($j = 1; $j -le $layoutShapes.Count; $j++) {
        $shp = $layoutShapes.Item($j)
        if ($shp.PlaceholderFormat.Type -eq 16) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}
